# "Generate Report for Handback": refresh the "Latest Handback DateTime"
# (column K) of the c71c5377-... row (row 2) on both locale report sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K2").Value = "2016-10-21 04:38:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-21 04:39:16"
